$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before the existing row 402, shifting the rest of
# the table (old rows 402-425) down to 404-427.
$ws.Rows("402:403").Insert()

# New record 1 (row 402): Larga vida / Primera, 02-06-2022 (serial 44714)
$ws.Range("A402").Value = 11
$ws.Range("B402").Value = "Vega Monumental Concepción"
$ws.Range("C402").Value = "Bíobío"
$ws.Range("D402").Value = 44714
$ws.Range("E402").Value = 8
$ws.Range("F402").Value = 100112020
$ws.Range("G402").Value = "Tomate"
$ws.Range("H402").Value = "Larga vida"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 850
$ws.Range("K402").Value = 14000
$ws.Range("L402").Value = 15000
$ws.Range("M402").Value = 14588
$ws.Range("N402").Value = "$/bandeja 18 kilos"
$ws.Range("O402").Value = "Región de Arica y Parinacota"
$ws.Range("P402").Value = 810
$ws.Range("Q402").Value = 18
$ws.Range("R402").Value = "Hortaliza"

# New record 2 (row 403): Larga vida / Segunda, 02-06-2022 (serial 44714)
$ws.Range("A403").Value = 11
$ws.Range("B403").Value = "Vega Monumental Concepción"
$ws.Range("C403").Value = "Bíobío"
$ws.Range("D403").Value = 44714
$ws.Range("E403").Value = 8
$ws.Range("F403").Value = 100112020
$ws.Range("G403").Value = "Tomate"
$ws.Range("H403").Value = "Larga vida"
$ws.Range("I403").Value = "Segunda"
$ws.Range("J403").Value = 600
$ws.Range("K403").Value = 11000
$ws.Range("L403").Value = 12000
$ws.Range("M403").Value = 11500
$ws.Range("N403").Value = "$/bandeja 18 kilos"
$ws.Range("O403").Value = "Región de Arica y Parinacota"
$ws.Range("P403").Value = 639
$ws.Range("Q403").Value = 18
$ws.Range("R403").Value = "Hortaliza"
